$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A33").Value = "GRT-USD"
